# Update countries & provincias Spain
# Refresh the "Pais" data table to the 21 Abril 2020 16:52 snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 16:52"

# --- Country name shifts (rows re-ranked by total cases) ------------------
# Crucero (row 92) and Albania (row 98) keep their spot; the countries in
# between shift down one rank and "Guinea" jumps up in front of "Libano".
$ws.Range("A93").Value = "Guinea"
$ws.Range("A94").Value = "Libano"
$ws.Range("A95").Value = "Nigeria"
$ws.Range("A96").Value = "Costa Rica"
$ws.Range("A97").Value = "Niger"

# "Siria" and "Mozambique" swap ranks.
$ws.Range("A169").Value = "Mozambique"
$ws.Range("A170").Value = "Siria"

# --- Updated case numbers ---------------------------------------------------
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 794322
$ws.Range("C4").Value = 1563
$ws.Range("E4").Value = 679348

# Row 8 - Reino Unido
$ws.Range("E8").Value = 47524
$ws.Range("G8").Value = 7
$ws.Range("H8").Value = 4869

# Row 64
$ws.Range("B64").Value = 1977
$ws.Range("C64").Value = 125
$ws.Range("E64").Value = 1012
$ws.Range("G64").Value = 470
$ws.Range("H64").Value = 489

# Row 68
$ws.Range("D68").Value = 357
$ws.Range("E68").Value = 1294

# Row 93 - now "Guinea"
$ws.Range("B93").Value = 688
$ws.Range("C93").Value = 66
$ws.Range("D93").Value = 127
$ws.Range("E93").Value = 555
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 6

# Row 94 - now "Libano"
$ws.Range("B94").Value = 677
$ws.Range("D94").Value = 103
$ws.Range("E94").Value = 553
$ws.Range("F94").Value = 27
$ws.Range("H94").Value = 21

# Row 95 - now "Nigeria"
$ws.Range("B95").Value = 665
$ws.Range("D95").Value = 188
$ws.Range("E95").Value = 455
$ws.Range("F95").Value = 2
$ws.Range("H95").Value = 22

# Row 96 - now "Costa Rica"
$ws.Range("B96").Value = 662
$ws.Range("C96").Value = 0
$ws.Range("E96").Value = 532
$ws.Range("F96").Value = 8
$ws.Range("H96").Value = 6

# Row 97 - now "Niger"
$ws.Range("B97").Value = 655
$ws.Range("C97").Value = 7
$ws.Range("D97").Value = 124
$ws.Range("E97").Value = 511
$ws.Range("H97").Value = 20

# Row 169 - now "Mozambique"
$ws.Range("D169").Value = 8
$ws.Range("H169").Value = 0

# Row 170 - now "Siria"
$ws.Range("D170").Value = 6
$ws.Range("E170").Value = 30
$ws.Range("H170").Value = 3
